# "Generate Report for Handoff"
#
# The handoff transform now succeeds for 7dda93ed-dfce-488e-8073-c70f64c3b907.md,
# so the status/report columns on the per-locale sheets move from
# "failed/ignored" to "ready/include", with the newly produced .xlf handoff
# files linked in and their handoff timestamps recorded.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text for the file, in both locale columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

# --- per-locale handoff report rows ---
$locales = @(
    @{ Sheet = "zh-cn"; Xlf = "7dda93ed-dfce-488e-8073-c70f64c3b907.01fa9251abc503a225750a4b743bc45ee34ed12a.zh-cn.xlf"; HandoffDate = "2016-01-25 08:20:32" },
    @{ Sheet = "de-de"; Xlf = "7dda93ed-dfce-488e-8073-c70f64c3b907.01fa9251abc503a225750a4b743bc45ee34ed12a.de-de.xlf"; HandoffDate = "2016-01-25 08:20:44" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Status (column B) - the handoff transform now succeeds
    $ws.Range("B2").Value = "Ready for handoff"

    # Latest Handoff File (column C) - link to the generated .xlf, styled like
    # the other hyperlink cells on this sheet (underlined, cornflower blue -
    # the workbook's "HyperLink" cell style)
    $ws.Range("C2").Value = $locale.Xlf
    $baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ba6004259b73dcc7b07181830d70be72dc883c71/e2e/"
    $ws.Hyperlinks.Add($ws.Range("C2"), $baseUrl + $locale.Xlf, "", "", $locale.Xlf)
    $ws.Range("C2").Style = "HyperLink"
    $ws.Range("C2").Font.Underline = 2
    $ws.Range("C2").Font.Color = 15570276

    # Latest Handoff Datetime (column D) - real timestamp instead of the zero date
    $ws.Range("D2").Value = $locale.HandoffDate

    # Handoff Reason (column H) - file is now included rather than ignored
    $ws.Range("H2").Value = "Include"
}
